$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts old row9 -> row10, etc.)
$ws.Rows.Item(9).Insert()

# Insert another new row at row 11 (after the old row9, now at row10)
$ws.Rows.Item(11).Insert()

# Fill new row 9: 7, 22, 223, 33, Plastic
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 22
$ws.Cells.Item(9, 3).Value = 223
$ws.Cells.Item(9, 4).Value = 33
$ws.Cells.Item(9, 5).Value = "Plastic"

# Fill new row 11: 9, 2, 2, 2, Wood
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(11, 5).Value = "Wood"
